# Re-order the "Recorded By" (column G) comma-separated list of recorders
# so that the literal entry "System" always appears first, with the
# remaining entries sorted alphabetically.
#
# Example: "dnasr281@gmail.com, System" -> "System, dnasr281@gmail.com"
#          "system, System, backup@backdoor.com" -> "System, backup@backdoor.com, system"
#          "dnasr281@gmail.com, admin@admin.com" -> "admin@admin.com, dnasr281@gmail.com"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

# Column G is the "Recorded By" column.
$col = 7

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $col)
    $value = $cell.Value()

    if ($null -eq $value) { continue }
    if (-not ($value -is [string])) { continue }
    if ($value.IndexOf(",") -lt 0) { continue }

    $parts = $value -split ","
    $trimmed = @()
    foreach ($p in $parts) {
        $trimmed += $p.Trim()
    }

    $systemParts = @()
    $otherParts = @()
    foreach ($p in $trimmed) {
        if ($p.Equals("System")) {
            $systemParts += $p
        } else {
            $otherParts += $p
        }
    }

    $otherSorted = $otherParts | Sort-Object

    $newParts = @()
    $newParts += $systemParts
    $newParts += $otherSorted

    $newValue = [string]::Join(", ", $newParts)

    if ($newValue -ne $value) {
        $cell.Value = $newValue
    }
}
